$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Name for row 2 (B2)
$ws.Range("B2").Value = "wdw"

# Update Active flags (column A) for rows 3, 7, 8, 9
$ws.Range("A3").Value = $true
$ws.Range("A7").Value = $false
$ws.Range("A8").Value = $false
$ws.Range("A9").Value = $true
